$p = $ppt.ActivePresentation

# This reverts a run-consolidation pass: PowerPoint had merged what used to be
# separate single-word/space runs (e.g. "Slide" + " " + "1" + " " + "(Content)")
# into fewer runs (e.g. "Slide " + "1 " + "(Content)"). We restore the original,
# more granular run boundaries by re-assigning the text of narrower
# sub-ranges (TextRange.Characters(start, length)) back onto themselves; the
# COM host splits the underlying <a:r> run whenever a Characters() sub-range
# assignment doesn't cover an entire existing run.

function Split-TitleRuns($slide, [string[]]$words) {
    $tr = $slide.Shapes.Item(1).TextFrame.TextRange
    $pos = 1
    foreach ($w in $words) {
        $tr.Characters($pos, $w.Length).Text = $w
        $pos += $w.Length + 1   # skip the single space that follows each word
    }
}

# Slides whose titles are "Slide N (Content)" / similar single-space-joined titles.
Split-TitleRuns $p.Slides.Item(1) @("Slide", "1")
Split-TitleRuns $p.Slides.Item(2) @("Slide", "2")
Split-TitleRuns $p.Slides.Item(3) @("Slide", "3")
Split-TitleRuns $p.Slides.Item(4) @("Slide", "4")
Split-TitleRuns $p.Slides.Item(5) @("Slide", "5", "(Two")
Split-TitleRuns $p.Slides.Item(6) @("Slide", "6", "(Two", "Content")
Split-TitleRuns $p.Slides.Item(7) @("Slide", "7", "(Content", "with")
Split-TitleRuns $p.Slides.Item(8) @("Slide", "8")
Split-TitleRuns $p.Slides.Item(9) @("Slide", "10")
Split-TitleRuns $p.Slides.Item(10) @("Slide", "11")
Split-TitleRuns $p.Slides.Item(11) @("Slide", "12")

# "an image" / "An image" captions that also had their leading word + space
# merged into a single run.
$p.Slides.Item(6).Shapes.Item(3).TextFrame.TextRange.Characters(1, 2).Text = "an"
$p.Slides.Item(7).Shapes.Item(4).TextFrame.TextRange.Characters(1, 2).Text = "An"
$p.Slides.Item(8).Shapes.Item(4).TextFrame.TextRange.Characters(1, 2).Text = "An"
